# Negative_Manifest sheet update
# - B6:B25 file names: IMG_XXXX index bumped by +1 (camera roll shifted by one)
# - rows 26:29 (n25-n28) removed entirely (last 4 negative samples dropped)
# - dimension shrinks from A1:E29 to A1:E25 accordingly
# - view: zoomed to 66%, selection re-anchored over A2:E29 -> A2:E25 range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the shifted file names in column B (rows 6-25) -----------------
$ws.Range("B6").Value  = "n5_IMG_3176.jpeg"
$ws.Range("B7").Value  = "n6_IMG_3176HorFlip.jpeg"
$ws.Range("B8").Value  = "n7_IMG_3176HorVertFlip.jpeg"
$ws.Range("B9").Value  = "n8_IMG_3176VertFlip.jpeg"
$ws.Range("B10").Value = "n9_IMG_3177.jpeg"
$ws.Range("B11").Value = "n10_IMG_3177HorFlip.jpeg"
$ws.Range("B12").Value = "n11_IMG_3177HorVertFlip.jpeg"
$ws.Range("B13").Value = "n12_IMG_3177VertFlip.jpeg"
$ws.Range("B14").Value = "n13_IMG_3178.jpeg"
$ws.Range("B15").Value = "n14_IMG_3178HorFlip.jpeg"
$ws.Range("B16").Value = "n15_IMG_3178HorVertFlip.jpeg"
$ws.Range("B17").Value = "n16_IMG_3178VertFlip.jpeg"
$ws.Range("B18").Value = "n17_IMG_3179.jpeg"
$ws.Range("B19").Value = "n18_IMG_3179HorFlip.jpeg"
$ws.Range("B20").Value = "n19_IMG_3179HorVertFlip.jpeg"
$ws.Range("B21").Value = "n20_IMG_3179VertFlip.jpeg"
$ws.Range("B22").Value = "n21_IMG_3180.jpeg"
$ws.Range("B23").Value = "n22_IMG_3180HorFlip.jpeg"
$ws.Range("B24").Value = "n23_IMG_3180HorVertFlip.jpeg"
$ws.Range("B25").Value = "n24_IMG_3180VertFlip.jpeg"

# --- Refresh view: zoom to 66% and (re)select the data range before the ----
# row count changes, matching the selection left behind in the saved file.
$ws.Range("A2:E29").Select() | Out-Null
$excel.ActiveWindow.Zoom = 66

# --- Remove the trailing rows (n25-n28), shrinking the manifest ------------
$ws.Rows("26:29").Delete()

$wb.Save()
